$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New "AVG" header next to the existing country columns on the second table (row 11)
$ws.Range("G11").Value = "AVG"

# Row 12 - "1 bedroom"
$ws.Range("B12").Value = 9266
$ws.Range("C12").Value = 4636
$ws.Range("D12").Value = 4523
$ws.Range("E12").Value = 32320
$ws.Range("F12").Value = 33597
$ws.Range("G12").Formula = "=(B12+C12+D12+E12+F12)/5"

# Row 13 - "2 bedrooms"
$ws.Range("B13").Value = 3487
$ws.Range("C13").Value = 1329
$ws.Range("D13").Value = 934
$ws.Range("E13").Value = 5570
$ws.Range("F13").Value = 8242
$ws.Range("G13").Formula = "=(B13+C13+D13+E13+F13)/5"

# Row 14 - "3 bedrooms"
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 402
$ws.Range("D14").Value = 382
$ws.Range("E14").Value = 1707
$ws.Range("F14").Value = 2185
$ws.Range("G14").Formula = "=(B14+C14+D14+E14+F14)/5"

# Row 15 - "4 bedrooms"
$ws.Range("B15").Value = 362
$ws.Range("C15").Value = 180
$ws.Range("D15").Value = 57
$ws.Range("E15").Value = 457
$ws.Range("F15").Value = 511
$ws.Range("G15").Formula = "=(B15+C15+D15+E15+F15)/5"

# Row 16 - "5+ bedrooms"
$ws.Range("B16").Formula = "=10+18+7+3+2"
$ws.Range("C16").Value = 49
$ws.Range("D16").Formula = "=14+4+4+1"
$ws.Range("E16").Formula = "=104+26+4+11+5+2"
$ws.Range("F16").Formula = "=3+6+4+23+95+1"
$ws.Range("G16").Formula = "=(B16+C16+D16+E16+F16)/5"

# Row 17 - new "Error ( 0 )" row
$ws.Range("A17").Value = "Error ( 0 )"
$ws.Range("B17").Value = 901
$ws.Range("C17").Value = 131
$ws.Range("D17").Value = 548
$ws.Range("E17").Value = 4037
$ws.Range("F17").Value = 11740
$ws.Range("G17").Formula = "=(B17+C17+D17+E17+F17)/5"

# Active selection moves to G21 per the saved view state
$null = $ws.Range("G21").Select()
